$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

# 1. Title: restructure from
#    "How Computational Tools and a Culture of Openness Supported the Collaboratively Authoring of a Text"
#    to
#    "Collaboratively Authoring a Text: The Roles of Computational Tools and a Culture of Openness"
$d.Content.Find.Execute(
    "How Computational Tools and a Culture of Openness Supported the Collaboratively Authoring of a Text",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Collaboratively Authoring a Text: The Roles of Computational Tools and a Culture of Openness",
    $wdReplaceAll)

# 2. Scholarly Significance paragraph - minor grammatical / wording changes

# 2a. "...highlighted parts of open science other than reproducibility..."
#     -> "...highlighted parts and stages of open science other than reproducibility..."
$d.Content.Find.Execute(
    "highlighted parts of open science other than reproducibility",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "highlighted parts and stages of open science other than reproducibility",
    $wdReplaceAll)

# 2b. " contributions from newcomers (and the public), and an overarching aim of making "
#     -> " contributions from newcomers, and an overarching aim of making "
$d.Content.Find.Execute(
    "contributions from newcomers (and the public), and an overarching aim of making",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "contributions from newcomers, and an overarching aim of making",
    $wdReplaceAll)

# 2c. " process and product widely-used. As this work aims to show, "
#     -> " process and product widely-used. As this project aims to show, "
$d.Content.Find.Execute(
    "process and product widely-used. As this work aims to show,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "process and product widely-used. As this project aims to show,",
    $wdReplaceAll)
